# Re-sort the "Periodo Mora" / "Valor Mora" columns (E:F) for the arrears
# table in ascending chronological order, for each employee block.
# Block 1: rows 16-22  (JORGE SANTIAGO PERTUZ MENDOZA, periods 1608-1702)
# Block 2: rows 23-109 (JOHN EDUARD ARREDONDO HUERTAS, periods 1710-2412)
# This mirrors the database update described in the commit message: the
# underlying data export changed the sort order of the period/arrears
# rows, which also corrects the pairing between some periods and their
# "Valor Mora" amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=16; Period="1608"; Valor=27580},
    @{Row=17; Period="1609"; Valor=27580},
    @{Row=18; Period="1610"; Valor=27580},
    @{Row=19; Period="1611"; Valor=27580},
    @{Row=20; Period="1612"; Valor=27580},
    @{Row=21; Period="1701"; Valor=27580},
    @{Row=22; Period="1702"; Valor=27580},
    @{Row=23; Period="1710"; Valor=7861},
    @{Row=24; Period="1711"; Valor=29480},
    @{Row=25; Period="1712"; Valor=29480},
    @{Row=26; Period="1801"; Valor=29480},
    @{Row=27; Period="1802"; Valor=29480},
    @{Row=28; Period="1803"; Valor=29480},
    @{Row=29; Period="1804"; Valor=29480},
    @{Row=30; Period="1805"; Valor=29480},
    @{Row=31; Period="1806"; Valor=29480},
    @{Row=32; Period="1807"; Valor=29480},
    @{Row=33; Period="1808"; Valor=29480},
    @{Row=34; Period="1809"; Valor=31249},
    @{Row=35; Period="1810"; Valor=31249},
    @{Row=36; Period="1811"; Valor=31249},
    @{Row=37; Period="1812"; Valor=31249},
    @{Row=38; Period="1901"; Valor=31249},
    @{Row=39; Period="1902"; Valor=31249},
    @{Row=40; Period="1903"; Valor=31249},
    @{Row=41; Period="1904"; Valor=31249},
    @{Row=42; Period="1905"; Valor=31249},
    @{Row=43; Period="1906"; Valor=31249},
    @{Row=44; Period="1907"; Valor=31249},
    @{Row=45; Period="1908"; Valor=31249},
    @{Row=46; Period="1909"; Valor=31249},
    @{Row=47; Period="1910"; Valor=31249},
    @{Row=48; Period="1911"; Valor=31249},
    @{Row=49; Period="1912"; Valor=31249},
    @{Row=50; Period="2001"; Valor=31249},
    @{Row=51; Period="2002"; Valor=31249},
    @{Row=52; Period="2003"; Valor=31249},
    @{Row=53; Period="2004"; Valor=31249},
    @{Row=54; Period="2005"; Valor=31249},
    @{Row=55; Period="2006"; Valor=31249},
    @{Row=56; Period="2007"; Valor=31249},
    @{Row=57; Period="2008"; Valor=31249},
    @{Row=58; Period="2009"; Valor=31249},
    @{Row=59; Period="2010"; Valor=31249},
    @{Row=60; Period="2011"; Valor=31249},
    @{Row=61; Period="2012"; Valor=31249},
    @{Row=62; Period="2101"; Valor=31249},
    @{Row=63; Period="2102"; Valor=31249},
    @{Row=64; Period="2103"; Valor=31249},
    @{Row=65; Period="2104"; Valor=31249},
    @{Row=66; Period="2105"; Valor=31249},
    @{Row=67; Period="2106"; Valor=31249},
    @{Row=68; Period="2107"; Valor=31249},
    @{Row=69; Period="2108"; Valor=31249},
    @{Row=70; Period="2109"; Valor=31249},
    @{Row=71; Period="2110"; Valor=31249},
    @{Row=72; Period="2111"; Valor=31249},
    @{Row=73; Period="2112"; Valor=31249},
    @{Row=74; Period="2201"; Valor=31249},
    @{Row=75; Period="2202"; Valor=31249},
    @{Row=76; Period="2203"; Valor=31249},
    @{Row=77; Period="2204"; Valor=31249},
    @{Row=78; Period="2205"; Valor=31249},
    @{Row=79; Period="2206"; Valor=31249},
    @{Row=80; Period="2207"; Valor=31249},
    @{Row=81; Period="2208"; Valor=31249},
    @{Row=82; Period="2209"; Valor=31249},
    @{Row=83; Period="2210"; Valor=31249},
    @{Row=84; Period="2211"; Valor=31249},
    @{Row=85; Period="2212"; Valor=31249},
    @{Row=86; Period="2301"; Valor=31249},
    @{Row=87; Period="2302"; Valor=31249},
    @{Row=88; Period="2303"; Valor=31249},
    @{Row=89; Period="2304"; Valor=31249},
    @{Row=90; Period="2305"; Valor=31249},
    @{Row=91; Period="2306"; Valor=31249},
    @{Row=92; Period="2307"; Valor=31249},
    @{Row=93; Period="2308"; Valor=31249},
    @{Row=94; Period="2309"; Valor=31249},
    @{Row=95; Period="2310"; Valor=31249},
    @{Row=96; Period="2311"; Valor=31249},
    @{Row=97; Period="2312"; Valor=31249},
    @{Row=98; Period="2401"; Valor=31249},
    @{Row=99; Period="2402"; Valor=31249},
    @{Row=100; Period="2403"; Valor=31249},
    @{Row=101; Period="2404"; Valor=31249},
    @{Row=102; Period="2405"; Valor=31249},
    @{Row=103; Period="2406"; Valor=31249},
    @{Row=104; Period="2407"; Valor=31249},
    @{Row=105; Period="2408"; Valor=31249},
    @{Row=106; Period="2409"; Valor=31249},
    @{Row=107; Period="2410"; Valor=31249},
    @{Row=108; Period="2411"; Valor=31249},
    @{Row=109; Period="2412"; Valor=20833}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Period
    $ws.Cells.Item($r.Row, 6).Value = $r.Valor
}
